$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = -0.04461569661326293
$ws.Range("B15").Value = -0.04461569661326292
$ws.Range("B16").Value = 0.00202798620969376
$ws.Range("B17").Value = 0.00202798620969376
$ws.Range("B18").Value = -0.0466436828229567
$ws.Range("B19").Value = -0.9999999999999976
$ws.Range("B23").Value = 0.9683634151287746
$ws.Range("B25").Value = 0.9237477185155117
$ws.Range("B26").Value = 0.9237477185155116
$ws.Range("B27").Value = 0.9237477185155117
$ws.Range("B28").Value = 1.608598661529098
$ws.Range("B33").Value = 0.05576962076657865
$ws.Range("B39").Value = 0.03285337659703792
$ws.Range("B40").Value = 0.06814033664570951
$ws.Range("B41").Value = -3.369701886027168
$ws.Range("B43").Value = 3.369701886027168
$ws.Range("B45").Value = 0.4295274792131411
$ws.Range("B47").Value = 0
$ws.Range("B48").Value = 0
$ws.Range("B49").Value = 0
$ws.Range("B100").Value = 3.422023930237267
$ws.Range("B104").Value = 2.171973230582027
$ws.Range("B118").Value = 0.4295274792131411
$ws.Range("B119").Value = 0.2147637396065706
$ws.Range("B120").Value = -3.344960454268904
$ws.Range("B121").Value = 3.344960454268904
$ws.Range("B122").Value = 0.9588318799432121
$ws.Range("B123").Value = -0.6546339484891495
$ws.Range("B124").Value = 0.6546339484891495
$ws.Range("B127").Value = 0.1423646319205026
$ws.Range("B128").Value = 0.1423646319205026
$ws.Range("B129").Value = 0
$ws.Range("B130").Value = 0.04502129385520167
$ws.Range("B133").Value = 0
$ws.Range("B136").Value = [double]"7.509768303535663e-17"
$ws.Range("B137").Value = 0
$ws.Range("B138").Value = -1.048874467653617
$ws.Range("B139").Value = 1.048874467653617
$ws.Range("B140").Value = -0.0255526262421415
$ws.Range("B141").Value = 0.02555262624214149
$ws.Range("B142").Value = 0.2261204623808553
$ws.Range("B143").Value = 0.5767592780369081
$ws.Range("B144").Value = -0.1046440884201997
$ws.Range("B145").Value = 0.1046440884201997
$ws.Range("B147").Value = 0.1216791725816262
$ws.Range("B148").Value = 0.1216791725816262
$ws.Range("B151").Value = 0.06814033664570951
$ws.Range("B175").Value = 0.8906915432975026
$ws.Range("B188").Value = 4.165686473331973
$ws.Range("B192").Value = -4.380450212938541
$ws.Range("B199").Value = -0.3930237274386525
$ws.Range("B202").Value = 0
$ws.Range("B203").Value = 0
$ws.Range("B204").Value = 0
$ws.Range("B205").Value = 0
$ws.Range("B208").Value = -3.422023930237267
$ws.Range("B210").Value = 0.9999999999999973
$ws.Range("B211").Value = 0
$ws.Range("B240").Value = 0.1423646319205026
$ws.Range("B241").Value = 0.04502129385520168
$ws.Range("B242").Value = 2.171973230582027
